# "The Scene And Player Readin Logic"
# Rescale the PlayerSpacing table's offset presets (divide the old
# pixel-ish magnitudes down to small gameplay units) and drop the
# now-unused 7th row (Key=7), shrinking the table from A1:C7 to A1:C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the preset value columns (Player_X / Player_Y) ------------
# Row 2 (Key=2) keeps "[0]" / "[0]" - unchanged.

# Row 3 (Key=3)
$ws.Range("B3").Value = "[-2,2]"
$ws.Range("C3").Value = "[0,0]"

# Row 4 (Key=4)
$ws.Range("B4").Value = "[-3,0,3]"
$ws.Range("C4").Value = "[1,1.5,1]"

# Row 5 (Key=5)
$ws.Range("B5").Value = "[-5,-2,2,5]"
$ws.Range("C5").Value = "[-0.5,2,2,-0.5]"

# Row 6 (Key=6)
$ws.Range("B6").Value = "[-5,-2.5,0,2.5,5]"
$ws.Range("C6").Value = "[-1,0,1,0,-1]"

# --- Drop the old row 7 (Key=7) -----------------------------------------
# This shifts the table/used-range up and resizes the table ref + sheet
# dimension from A1:C7 down to A1:C6 automatically.
$ws.Rows.Item(7).Delete() | Out-Null

# --- Match the saved selection state ------------------------------------
$ws.Range("C9").Activate() | Out-Null
$ws.Range("C9:C10").Select() | Out-Null
